# Adds three new worksheets (Baseline-South-CA, Baseline-Central-CA,
# Expanded-all-ports) after the existing "Baseline-limited-ports" sheet,
# each populated with a Year / Cummulative Capacity table that mirrors
# the layout and header styling of the original sheet.

$wb = $excel.ActiveWorkbook
$srcSheet = $wb.Worksheets.Item(1)

function Add-CapacitySheet {
    # NOTE: named parameter binding (-SheetName ...) is unreliable in this
    # engine, so this function is always called with positional arguments:
    # Add-CapacitySheet <name> <years[]> <values[]>
    param(
        [string]$SheetName,
        [object[]]$Years,
        [object[]]$Values
    )

    $lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
    $ws = $wb.Worksheets.Add([Type]::Missing, $lastSheet)
    $ws.Name = $SheetName

    # Header row
    $ws.Range("A1").Value = "Year"
    $ws.Range("B1").Value = "Cummulative Capacity"

    # Copy header formatting (bold, borders, centered/top aligned) from the
    # source sheet's header cells so the new sheets match the original style.
    $srcSheet.Range("A1:B1").Copy()
    $ws.Range("A1:B1").PasteSpecial(-4122)

    # Data rows
    $n = $Years.Length
    $data = New-Object 'object[,]' $n,2
    for ($i = 0; $i -lt $n; $i++) {
        $data[$i,0] = $Years[$i]
        $data[$i,1] = $Values[$i]
    }
    $lastRow = 1 + $n
    $ws.Range("A2:B$lastRow").Value = $data

    return $ws
}

# --- Baseline-South-CA ---
$years1 = 2030,2031,2032,2033,2034,2035,2036,2037,2038,2039,2040,2041,2042,2043,2044,2045
$vals1 = 407.2304176189487,1407.230417618949,2416.371618750885,4046.408984503963,6266.177250188815,9220.979273168934,12105.86042738475,15304.58631049695,17932.10799615077,21257.82231210306,24375.96813480699,27564.34172178803,29543.23725055432,31079.73360091929,32543.23725055432,33898.03721359599
Add-CapacitySheet "Baseline-South-CA" $years1 $vals1 | Out-Null

# --- Baseline-Central-CA ---
$years2 = 2030,2031,2032,2033,2034,2035,2036,2037,2038,2039,2040,2041,2042,2043,2044,2045,2046,2047,2048,2049
$vals2 = 407.2304176189487,1407.230417618949,2416.371618750885,3445.629284073075,5066.24564964162,6822.730240215845,8443.038978354918,10239.4110814471,12665.8387943053,15919.61597541799,19173.43735505323,22426.13538510296,24540.70647080056,26139.88567736691,27804.48371464776,29960.8166642758,31862.12245642851,32661.46546217721,33460.76642335767,34262.2502053107
Add-CapacitySheet "Baseline-Central-CA" $years2 $vals2 | Out-Null

# --- Expanded-all-ports ---
$years3 = 2030,2031,2032,2033,2034,2035,2036,2037,2038,2039,2040,2041,2042,2043,2044,2045,2046,2047,2048,2049
$vals3 = 407.2304176189487,1407.230417618949,2416.371618750885,4046.408984503963,6266.177250188815,8621.841047329388,10841.32899203481,13238.52188856064,15664.99058214565,18984.76240666921,23522.17776208614,28601.34058502495,33669.33764385931,38737.44076710031,43795.9591328035,48604.10730007126,51764.77980879939,53198.10201660736,53997.44502235606,54808.9430894309
Add-CapacitySheet "Expanded-all-ports" $years3 $vals3 | Out-Null

# Restore selection/active sheet to the first sheet, matching original workbook state.
$srcSheet.Select() | Out-Null
$srcSheet.Range("A1").Select() | Out-Null

Write-Host "Added sheets. Total worksheets: $($wb.Worksheets.Count)"
